$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.714.02'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.633.41'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '212.06'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '23.20'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').Value = '0.0612'
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('E11').Value = '  -3.12%  '
$ws.Range('D12').Value = '1.865.62'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '1.630.49'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '0.555'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').Value = '65.19'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '27.671.10'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '230.02'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').Value = '7.56'
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('E22').Value = '  +4.57%  '
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('D25').Value = '148.78'
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('E26').Value = '  -1.31%  '
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').Value = '15.60'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('D32').Value = '3.28'
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').Value = '1.471.97'
$ws.Range('E33').Value = '  +0.45%  '
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('E37').Value = '  +5.31%  '
$ws.Range('D38').Value = '0.877'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = '0.558'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').Value = '67.92'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D43').Value = '2.48'
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('E44').Value = '  -1.08%  '
$ws.Range('E45').Value = '  -4.60%  '
$ws.Range('D46').Value = '1.774.54'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').Value = '87.58'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').Value = '0.0994'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').Value = '7.72'
$ws.Range('E51').Value = '  -1.32%  '
